$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (text corrections for Sarah Wald profile) ---

# Organization name + URL corrected to Cambridge University
$ws.Range("A27").Value = "Cambridge University"
$ws.Range("B27").Value = "www.cam.ac.uk"

# Instagram handle -> full www. URL
$ws.Range("B5").Value = "www.instagram.com/sarahwald_/"

# Role: add affiliation
$ws.Range("B4").Value = "Summer intern, University of Cambridge"

# Groups: corrected group name
$ws.Range("B15").Value = "Visiting and Co-supervised Students"

# --- Add hyperlinks for the corrected URLs ---
# (rId3 -> B27, rId4 -> B5, matching creation order)
$ws.Hyperlinks.Add($ws.Range("B27"), "http://www.cam.ac.uk/")
$ws.Hyperlinks.Add($ws.Range("B5"), "http://www.instagram.com/sarahwald_/")

# --- Adjust the view back to a "normal" scroll/selection state ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B15").Select()
